$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E).
# Numeric-looking Price values are written with a leading apostrophe so Excel
# keeps them as text (matching the original inline-string cell content) instead
# of auto-converting them to floating-point numbers.

$ws.Range("D2").Value = "41.151.39"
$ws.Range("E2").Value = "  -6.31%  "
$ws.Range("D3").Value = "2.210.32"
$ws.Range("E3").Value = "  -6.69%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'242.77"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("E6").Value = "  -7.33%  "
$ws.Range("D7").Value = "'69.85"
$ws.Range("E7").Value = "  -5.95%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.551"
$ws.Range("E9").Value = "  -9.38%  "
$ws.Range("D10").Value = "'37.31"
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("D11").Value = "'0.0949"
$ws.Range("E11").Value = "  -7.60%  "
$ws.Range("D12").Value = "'57.73"
$ws.Range("E12").Value = "  -3.27%  "
$ws.Range("E13").Value = "  -4.25%  "
$ws.Range("E14").Value = "  -8.63%  "
$ws.Range("D15").Value = "2.540.55"
$ws.Range("E15").Value = "  -6.88%  "
$ws.Range("E16").Value = "  -10.04%  "
$ws.Range("D17").Value = "'0.837"
$ws.Range("E17").Value = "  -9.86%  "
$ws.Range("D18").Value = "2.209.81"
$ws.Range("E18").Value = "  -7.04%  "
$ws.Range("D19").Value = "41.138.53"
$ws.Range("E19").Value = "  -6.41%  "
$ws.Range("E20").Value = "  -8.43%  "
$ws.Range("D21").Value = "'72.66"
$ws.Range("E21").Value = "  -6.75%  "
$ws.Range("E22").Value = "  -8.46%  "
$ws.Range("D23").Value = "'231.35"
$ws.Range("E23").Value = "  -9.10%  "
$ws.Range("D24").Value = "'2.01"
$ws.Range("E24").Value = "  +6.43%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -5.28%  "
$ws.Range("E27").Value = "  -3.85%  "
$ws.Range("E28").Value = "  -5.12%  "
$ws.Range("E29").Value = "  -8.04%  "
$ws.Range("D30").Value = "'171.60"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("D31").Value = "'20.38"
$ws.Range("E31").Value = "  -9.18%  "
$ws.Range("D32").Value = "'0.118"
$ws.Range("E32").Value = "  -9.13%  "
$ws.Range("E33").Value = "  -8.07%  "
$ws.Range("D34").Value = "'0.0706"
$ws.Range("E34").Value = "  -7.33%  "
$ws.Range("D35").Value = "'5.15"
$ws.Range("E35").Value = "  -5.22%  "
$ws.Range("E36").Value = "  -9.93%  "
$ws.Range("D37").Value = "'3.88"
$ws.Range("E37").Value = "  +1.63%  "
$ws.Range("D38").Value = "'23.71"
$ws.Range("E38").Value = "  +14.75%  "
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("E40").Value = "  -5.92%  "
$ws.Range("D41").Value = "'5.81"
$ws.Range("E41").Value = "  -12.52%  "
$ws.Range("D42").Value = "'63.76"
$ws.Range("E42").Value = "  -3.35%  "
$ws.Range("D43").Value = "'4.95"
$ws.Range("E43").Value = "  -10.28%  "
$ws.Range("E44").Value = "  -4.67%  "
$ws.Range("D45").Value = "'8.60"
$ws.Range("E45").Value = "  -5.45%  "
$ws.Range("D46").Value = "'0.0999"
$ws.Range("E46").Value = "  -7.19%  "
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("E48").Value = "  +9.07%  "
$ws.Range("D49").Value = "'4.43"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").Value = "'1.16"
$ws.Range("E50").Value = "  -6.96%  "
$ws.Range("E51").Value = "  -6.06%  "
